# Auto-generated: applies the crypto-price refresh described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '79.674.81'
$ws.Range('E2').Value = '  +4.35%  '
$ws.Range('D3').Value = '3.199.54'
$ws.Range('E3').Value = '  +5.33%  '
$ws.Range('D5').Value = "'" + '205.84'
$ws.Range('E5').Value = '  +2.69%  '
$ws.Range('D6').Value = "'" + '635.87'
$ws.Range('E6').Value = '  +1.74%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = "'" + '0.245'
$ws.Range('E8').Value = '  +19.02%  '
$ws.Range('D9').Value = "'" + '0.611'
$ws.Range('E9').Value = '  +11.33%  '
$ws.Range('D10').Value = '3.198.57'
$ws.Range('E10').Value = '  +5.35%  '
$ws.Range('D11').Value = "'" + '0.629'
$ws.Range('E11').Value = '  +43.29%  '
$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D12').Value = "'" + '0.0000245'
$ws.Range('E12').Value = '  +26.96%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = "'" + '0.166'
$ws.Range('E13').Value = '  +3.31%  '
$ws.Range('E14').Value = '  +4.18%  '
$ws.Range('D15').Value = '3.785.48'
$ws.Range('E15').Value = '  +5.23%  '
$ws.Range('D16').Value = "'" + '32.50'
$ws.Range('E16').Value = '  +11.65%  '
$ws.Range('D17').Value = '79.448.27'
$ws.Range('E17').Value = '  +4.04%  '
$ws.Range('D18').Value = '3.195.41'
$ws.Range('E18').Value = '  +5.28%  '
$ws.Range('D19').Value = "'" + '14.68'
$ws.Range('E19').Value = '  +9.11%  '
$ws.Range('D20').Value = "'" + '9.40'
$ws.Range('E20').Value = '  +4.68%  '
$ws.Range('D21').Value = "'" + '2.98'
$ws.Range('E21').Value = '  +29.72%  '
$ws.Range('D22').Value = "'" + '435.36'
$ws.Range('E22').Value = '  +16.22%  '
$ws.Range('D23').Value = "'" + '5.24'
$ws.Range('E23').Value = '  +20.21%  '
$ws.Range('D24').Value = "'" + '4.88'
$ws.Range('E24').Value = '  +11.59%  '
$ws.Range('D25').Value = '3.366.09'
$ws.Range('E25').Value = '  +5.50%  '
$ws.Range('D26').Value = "'" + '77.73'
$ws.Range('E26').Value = '  +6.22%  '
$ws.Range('D27').Value = "'" + '11.08'
$ws.Range('E27').Value = '  +12.87%  '
$ws.Range('D28').Value = "'" + '1.00'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = "'" + '0.0000120'
$ws.Range('E29').Value = '  +9.50%  '
$ws.Range('D30').Value = "'" + '9.30'
$ws.Range('E30').Value = '  +12.62%  '
$ws.Range('D31').Value = "'" + '0.999'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  +6.41%  '
$ws.Range('D33').Value = "'" + '530.65'
$ws.Range('E33').Value = '  +7.61%  '
$ws.Range('E34').Value = '  +3.05%  '
$ws.Range('E35').Value = '  +23.84%  '
$ws.Range('D36').Value = "'" + '23.27'
$ws.Range('E36').Value = '  +12.78%  '
$ws.Range('D37').Value = "'" + '0.125'
$ws.Range('E37').Value = '  +19.04%  '
$ws.Range('D39').Value = "'" + '0.412'
$ws.Range('E39').Value = '  +7.16%  '
$ws.Range('D40').Value = "'" + '164.63'
$ws.Range('E40').Value = '  +1.13%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').Value = "'" + '192.35'
$ws.Range('E42').Value = '  +1.48%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').Value = "'" + '5.55'
$ws.Range('E44').Value = '  +8.26%  '
$ws.Range('E45').Value = '  +10.08%  '
$ws.Range('D46').Value = "'" + '0.805'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('E47').Value = '  +5.52%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = "'" + '43.62'
$ws.Range('E48').Value = '  +3.78%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = "'" + '2.60'
$ws.Range('E49').Value = '  +5.79%  '
$ws.Range('D50').Value = "'" + '25.82'
$ws.Range('E50').Value = '  +15.69%  '
$ws.Range('D51').Value = "'" + '0.640'
$ws.Range('E51').Value = '  +5.41%  '
